$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 17 and 18, pushing the existing rows 17-28 down to 19-30.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(18).Insert()

# Populate the new row 17 with its data.
$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 45118
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100112013
$ws.Cells.Item(17, 7).Value = "Alcachofa"
$ws.Cells.Item(17, 8).Value = "Española"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 360
$ws.Cells.Item(17, 11).Value = 14000
$ws.Cells.Item(17, 12).Value = 15000
$ws.Cells.Item(17, 13).Value = 14500
$ws.Cells.Item(17, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(17, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 16).Value = 483
$ws.Cells.Item(17, 17).Value = 30
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Populate the new row 18 with its data.
$ws.Cells.Item(18, 1).Value = 8
$ws.Cells.Item(18, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 45118
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = 100112013
$ws.Cells.Item(18, 7).Value = "Alcachofa"
$ws.Cells.Item(18, 8).Value = "Madrigal"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 400
$ws.Cells.Item(18, 11).Value = 12500
$ws.Cells.Item(18, 12).Value = 13000
$ws.Cells.Item(18, 13).Value = 12750
$ws.Cells.Item(18, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(18, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(18, 16).Value = 319
$ws.Cells.Item(18, 17).Value = 40
$ws.Cells.Item(18, 18).Value = "Hortaliza"
